$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:F) to (B:G).
# This also carries column A's old "row label" formatting (bold/centered/bordered)
# into the new column B.
$ws.Columns("A:A").Insert()

# New header for the inserted column, in row 1.
$ws.Range("B1").Value = "segments"

# Give the new header cell the same header formatting as the rest of row 1.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill the new column A (rows 2-20) with the zero-based segment index.
for ($i = 0; $i -le 18; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

# Column A should carry the "row label" formatting that used to live on the
# segment-name column (now shifted into column B).
$ws.Range("B2:B20").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The segment-name text in column B reverts to plain (unstyled) formatting,
# matching how the numeric data columns are formatted.
$ws.Range("B2:B20").ClearFormats()
